# [UPDATE] can import types and read meta sheet
#
# The "Meta" sheet's header row described the second column as
# "Assembly Qualified Type Name" - rename it to the shorter "Type"
# (the cell below it keeps its existing value, the fully qualified
# type name used for import).
#
# Also restore the active sheet/selection so the workbook opens on the
# "Meta" sheet (cell B1 selected) instead of "Resources".

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("Meta")
$resourcesSheet = $wb.Worksheets.Item("Resources")

# Rename the column header used to describe the resource item's type.
$metaSheet.Range("B1").Value = "Type"

# Make "Meta" the active/selected sheet, with B1 selected - "Resources"
# goes back to being a plain, non-active tab (keeping its own selection).
$metaSheet.Activate()
$metaSheet.Range("B1").Select()
